$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = 42587.832430555558

$ws.Range("B5").Value = "Bag"
$ws.Range("C5").Value = 21
$ws.Range("D5").Value = 17
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 100
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = 0
